# Update CDL_NVC_AgClassMatch.xlsx: address wheat/corn issue
# For rows 3-53 (excluding rows that already have a value in column D),
# set column D to "Wheat". Row 30 already had a (misplaced) value in D,
# so that existing value is shifted right into column E first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 is special: it already has a value in D30 ("Row Crop - Close Grown
# Crop") that needs to move over to E30 before D30 becomes "Wheat".
$existing = $ws.Range("D30").Value()
$ws.Range("E30").Value = $existing
$ws.Range("D30").Value = "Wheat"

# Remaining rows just need "Wheat" written into the empty D cell.
$rows = @(3,4,5,6,7,8,9,10,11,12,13,22,24,25,26,27,28,29,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,49,50,51,52,53)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "Wheat"
}

# Update the view to reflect where the author ended up scrolled/selected.
$ws.Activate()
$ws.Range("A34").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F50").Select()
